# "språkvask" (language proofreading) pass on the QField tree-register app
# attribute table: fix a typo in a repeated shared-string label, and carry
# over the author's final cursor position / scroll state in the sheet view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Typo fix: "VAT19 variabler" -> "VAT variabler"
#    This label is repeated for every row of the VAT-variable block,
#    B45:B61, all pointing at the same shared string.
# ---------------------------------------------------------------------
$ws.Range("B45:B61").Value = "VAT variabler"

# ---------------------------------------------------------------------
# 2) View state: scroll so row 10 is at the top and leave the cursor on B17.
# ---------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B17").Select()

# ---------------------------------------------------------------------
# 3) Row heights for the re-flowed block (rows 45-61) settle a touch
#    shorter once the text is re-measured.
# ---------------------------------------------------------------------
$ws.Rows("45:61").RowHeight = 13.8
